$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the extra, empty "Sheet1" tab that shipped with the workbook.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

# ---------------------------------------------------------------------------
# 2. Turn the data range on "+informal-country" into an Excel Table (ListObject).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("+informal-country")
$rng = $ws.Range("A1:L165")
$lo = $ws.ListObjects.Add(1, $rng, 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"

# Column names (already present as the header row, but make sure the
# ListColumns line up with the expected names explicitly).
$names = @("country","GRIP-region","IMAGE-region","Impassable","Very horrible","Horrible","Bad","Intermediate","Good","Excellent","Total","Weighted grade")
for ($i = 0; $i -lt $names.Length; $i++) {
    $lo.ListColumns.Item($i + 1).Name = $names[$i]
}

# ---------------------------------------------------------------------------
# 3. Apply the "0.00" numeric display format to the Total / Weighted grade
#    columns (mirrors the two new dxf entries in styles.xml).
# ---------------------------------------------------------------------------
$lo.ListColumns.Item("Total").DataBodyRange.NumberFormat = "0.00"
$lo.ListColumns.Item("Weighted grade").DataBodyRange.NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 4. Cosmetic sheet-view tweaks: column widths + active selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.0
$ws.Columns.Item(2).ColumnWidth = 13.0
$ws.Columns.Item(3).ColumnWidth = 14.83
$ws.Columns.Item(4).ColumnWidth = 12.33
$ws.Columns.Item(5).ColumnWidth = 14.0
$ws.Columns.Item(6).ColumnWidth = 9.67
$ws.Columns.Item(8).ColumnWidth = 14.0
$ws.Columns.Item(10).ColumnWidth = 10.5
$ws.Columns.Item(12).ColumnWidth = 16.67

$ws.Activate()
$ws.Range("E5").Select()

Write-Host "done"
